$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 31 (pushing the
# existing rows 31-69 down to 32-70), matching how this "diario" sheet is
# kept as a rolling weekly log with the newest entry inserted near the top.
$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44540
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100112052
$ws.Range("G31").Value = "Albahaca"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = 3500
$ws.Range("N31").Value = "$/paquete"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 3500
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
